# case : update database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 from the company name string "FIRSTSTEP" to the numeric value 111
$ws.Range("B2").Value = 111

# Update the active selection on the sheet from D13 to B8
$ws.Range("B8").Select()
